$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D. This shifts the existing
#    period columns D:K (old period columns) one column to the right,
#    to E:L, making room for a new (most recent) reporting period in D.
$ws.Columns("D").Insert()

# 2. The newly inserted column D has no formatting of its own yet;
#    copy the number formats/styles from column E (which now holds
#    what used to be in D) so the new column visually matches its
#    neighbours (date format for the header row, number format for
#    data rows, etc).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# 3. Populate the new column D with the new period's figures.
#    Row 7/38/80 = "Period Ending" headers (new date 2018-12-31 = 43465)
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(80, 4).Value = 43465

# Income Statement (rows 8-35)
$ws.Cells.Item(8, 4).Value = 137100
$ws.Cells.Item(9, 4).Value = "NA"
$ws.Cells.Item(10, 4).Value = "NA"
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(15, 4).Value = "NA"
$ws.Cells.Item(17, 4).Value = 45800
$ws.Cells.Item(18, 4).Value = 91200
$ws.Cells.Item(20, 4).Value = -53400
$ws.Cells.Item(21, 4).Value = 38700
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(23, 4).Value = 37800
$ws.Cells.Item(24, 4).Value = 1600
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 36200
$ws.Cells.Item(27, 4).Value = 36200
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = 53400
$ws.Cells.Item(33, 4).Value = 36200
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 36200

# Balance Sheet (rows 41-77)
$ws.Cells.Item(41, 4).Value = 89500
$ws.Cells.Item(42, 4).Value = 270900
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 20100
$ws.Cells.Item(49, 4).Value = 3900
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 20300
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 5163900
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(59, 4).Value = 0
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(61, 4).Value = 238500
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 4863500
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 301500
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 300400
$ws.Cells.Item(77, 4).Value = 0

# Cash Flow Statement (rows 81-102)
$ws.Cells.Item(81, 4).Value = 36200
$ws.Cells.Item(83, 4).Value = 900
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 40800
$ws.Cells.Item(91, 4).Value = -3600
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -392400
$ws.Cells.Item(96, 4).Value = -2200
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = 337600
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(102, 4).Value = -13900

# 4. Two rows also carry a small restatement of the *prior* (now in
#    column E) period's figure alongside the new period's number, so
#    column E needs an explicit correction rather than simply being
#    the shifted-over old column D value.
$ws.Cells.Item(89, 5).Value = 21600
$ws.Cells.Item(94, 5).Value = -194300
